$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 247.625
$ws.Range("I12").Value = 221.83333
$ws.Range("J12").Value = 325
$ws.Range("K12").Value = 221.83333
$ws.Range("L12").Value = 325
$ws.Range("M12").Value = -51.83332999999999
$ws.Range("N12").Value = -665

$ws.Range("H40").Value = 2990
$ws.Range("I40").Value = 2800
$ws.Range("J40").Value = 3433.3333
$ws.Range("K40").Value = 2800
$ws.Range("L40").Value = 3433.3333
$ws.Range("M40").Value = -2625
$ws.Range("N40").Value = -3783.3333

$ws.Range("H62").Value = 35717210
$ws.Range("I62").Value = 41669508
$ws.Range("J62").Value = 3400
$ws.Range("K62").Value = 41669508
$ws.Range("L62").Value = 3400
$ws.Range("M62").Value = -41668884
$ws.Range("N62").Value = -4648

$ws.Range("H65").Value = 35717210
$ws.Range("I65").Value = 41669508
$ws.Range("J65").Value = 3400
$ws.Range("K65").Value = 208347540
$ws.Range("L65").Value = 17000
$ws.Range("M65").Value = -208344420
$ws.Range("N65").Value = -23240

$ws.Range("H113").Value = 1902.6666
$ws.Range("I113").Value = 1771.5
$ws.Range("J113").Value = 1957.8948
$ws.Range("K113").Value = 1771.5
$ws.Range("L113").Value = 1957.8948
$ws.Range("M113").Value = 1482.5
$ws.Range("N113").Value = -8465.8948

$ws.Range("H132").Value = 1910.12
$ws.Range("I132").Value = 1910.12
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 5730.36
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -3200.36

$ws.Range("H137").Value = 13160216
$ws.Range("I137").Value = 25001790
$ws.Range("J137").Value = 2910.889
$ws.Range("K137").Value = 75005370
$ws.Range("L137").Value = 8732.667000000001
$ws.Range("M137").Value = -75002820
$ws.Range("N137").Value = -13832.667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2571.2307
$ws.Range("I2").Value = 2719.3635
$ws.Range("J2").Value = 1756.5
$ws.Range("K2").Value = 2719.3635
$ws.Range("L2").Value = 1756.5
$ws.Range("M2").Value = -2606.3635
$ws.Range("N2").Value = -1982.5

$ws.Range("H45").Value = 1525.875
$ws.Range("I45").Value = 1290.591
$ws.Range("J45").Value = 4114
$ws.Range("K45").Value = 1290.591
$ws.Range("L45").Value = 4114
$ws.Range("M45").Value = -913.5909999999999
$ws.Range("N45").Value = -4868

$ws.Range("H61").Value = 3686.4443
$ws.Range("I61").Value = 1161.8334
$ws.Range("J61").Value = 8735.666999999999
$ws.Range("K61").Value = 1161.8334
$ws.Range("L61").Value = 8735.666999999999
$ws.Range("M61").Value = -949.8334
$ws.Range("N61").Value = -9159.666999999999

$ws.Range("H63").Value = 3995.1428
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 3995.1428
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 3995.1428
$ws.Range("N63").Value = -5367.1428

$ws.Range("H66").Value = 3995.1428
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 3995.1428
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 19975.714
$ws.Range("N66").Value = -26839.714

$ws.Range("H74").Value = 11898.7
$ws.Range("I74").Value = 12552.556
$ws.Range("J74").Value = 6014
$ws.Range("K74").Value = 12552.556
$ws.Range("L74").Value = 6014
$ws.Range("M74").Value = -11678.556
$ws.Range("N74").Value = -7762

$ws.Range("H77").Value = 11898.7
$ws.Range("I77").Value = 12552.556
$ws.Range("J77").Value = 6014
$ws.Range("K77").Value = 62762.78
$ws.Range("L77").Value = 30070
$ws.Range("M77").Value = -58394.78
$ws.Range("N77").Value = -38806

$ws.Range("H110").Value = 996.8461
$ws.Range("I110").Value = 746.7619
$ws.Range("J110").Value = 2047.2
$ws.Range("K110").Value = 746.7619
$ws.Range("L110").Value = 2047.2
$ws.Range("M110").Value = 1298.2381
$ws.Range("N110").Value = -6137.2

$ws.Range("H116").Value = 2571.2307
$ws.Range("I116").Value = 2719.3635
$ws.Range("J116").Value = 1756.5
$ws.Range("K116").Value = 2719.3635
$ws.Range("L116").Value = 1756.5
$ws.Range("M116").Value = -425.3634999999999
$ws.Range("N116").Value = -6344.5

$ws.Range("H136").Value = 3686.4443
$ws.Range("I136").Value = 1161.8334
$ws.Range("J136").Value = 8735.666999999999
$ws.Range("K136").Value = 3485.5002
$ws.Range("L136").Value = 26207.001
$ws.Range("M136").Value = -935.5001999999999
$ws.Range("N136").Value = -31307.001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2571.2307
$ws.Range("I3").Value = 2719.3635
$ws.Range("J3").Value = 1756.5
$ws.Range("K3").Value = 2719.3635
$ws.Range("L3").Value = 1756.5
$ws.Range("M3").Value = -2605.3635
$ws.Range("N3").Value = -1984.5

$ws.Range("H94").Value = 379.0435
$ws.Range("I94").Value = 375.44446
$ws.Range("J94").Value = 392
$ws.Range("K94").Value = 375.44446
$ws.Range("L94").Value = 392
$ws.Range("M94").Value = 75.55554000000001
$ws.Range("N94").Value = -1294

$ws.Range("H105").Value = 2787.1538
$ws.Range("I105").Value = 1347.5
$ws.Range("J105").Value = 3005.8354
$ws.Range("K105").Value = 1347.5
$ws.Range("L105").Value = 3005.8354
$ws.Range("M105").Value = 399.5
$ws.Range("N105").Value = -6499.8354

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2315.3242
$ws.Range("I31").Value = 1678.4231
$ws.Range("J31").Value = 3820.7273
$ws.Range("K31").Value = 1678.4231
$ws.Range("L31").Value = 3820.7273
$ws.Range("M31").Value = -1383.4231
$ws.Range("N31").Value = -4410.7273

$ws.Range("H34").Value = 2315.3242
$ws.Range("I34").Value = 1678.4231
$ws.Range("J34").Value = 3820.7273
$ws.Range("K34").Value = 1678.4231
$ws.Range("L34").Value = 3820.7273
$ws.Range("M34").Value = -1476.4231
$ws.Range("N34").Value = -4224.7273

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 14000
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 14000
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 14000
$ws.Range("N26").Value = -14560

$ws.Range("H50").Value = 14000
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 14000
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 14000
$ws.Range("N50").Value = -14996

$ws.Range("H122").Value = 1521.6
$ws.Range("I122").Value = 1350
$ws.Range("J122").Value = 2208
$ws.Range("K122").Value = 4050
$ws.Range("L122").Value = 6624
$ws.Range("M122").Value = -1600
$ws.Range("N122").Value = -11524

$ws.Range("H126").Value = 3189.3823
$ws.Range("I126").Value = 2834.7083
$ws.Range("J126").Value = 4040.6
$ws.Range("K126").Value = 8504.124899999999
$ws.Range("L126").Value = 12121.8
$ws.Range("M126").Value = -6034.124899999999
$ws.Range("N126").Value = -17061.8

$ws.Range("H132").Value = 5996
$ws.Range("I132").Value = 6394.8
$ws.Range("J132").Value = 5198.4
$ws.Range("K132").Value = 19184.4
$ws.Range("L132").Value = 15595.2
$ws.Range("M132").Value = -16654.4
$ws.Range("N132").Value = -20655.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2838.125
$ws.Range("I7").Value = 2916.6667
$ws.Range("J7").Value = 2602.5
$ws.Range("K7").Value = 2916.6667
$ws.Range("L7").Value = 2602.5
$ws.Range("M7").Value = -2804.6667
$ws.Range("N7").Value = -2826.5

$ws.Range("H122").Value = 2845.4546
$ws.Range("I122").Value = 2975
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 8925
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -6475
$ws.Range("N122").Value = -12400

$ws.Range("H126").Value = 2838.125
$ws.Range("I126").Value = 2916.6667
$ws.Range("J126").Value = 2602.5
$ws.Range("K126").Value = 8750.000100000001
$ws.Range("L126").Value = 7807.5
$ws.Range("M126").Value = -6280.000100000001
$ws.Range("N126").Value = -12747.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2303.147
$ws.Range("I132").Value = 1796.2593
$ws.Range("J132").Value = 4258.2856
$ws.Range("K132").Value = 5388.7779
$ws.Range("L132").Value = 12774.8568
$ws.Range("M132").Value = -2858.7779
$ws.Range("N132").Value = -17834.8568
